# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig
# Rebrand from "ibm.com" / "Alvearie Team" to "linuxforhealth.org" / "LinuxForHealth Team",
# bump the StructureDefinition version/date.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-detail"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

$elem.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-compared-to}
"
$elem.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-criteria}
"
$elem.Range("J7").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-level}
"
$elem.Range("J8").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-score}
"

# Q9 shares the same string as the Metadata URL cell (http://ibm.com/...match-detail),
# so it must be rebranded too.
$elem.Range("Q9").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-detail"
